$d = $word.ActiveDocument

$wNs  = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$w14Ns = "http://schemas.microsoft.com/office/word/2010/wordml"

# ---------------------------------------------------------------------------
# Hunk 1: paragraph "Jonathan Matetzky" (w14:paraId="1BF51220") gains RTL
# paragraph-mark run properties (rFonts hint=cs, rtl) in its pPr/rPr.
# ---------------------------------------------------------------------------
$target1 = '<w:p xmlns:w="' + $wNs + '" xmlns:w14="' + $w14Ns + '" w14:paraId="1BF51220" w14:textId="701F9F56" w:rsidR="00C87AA4" w:rsidRDefault="00C87AA4" w:rsidP="00C1464E">' + `
  '<w:pPr>' + `
    '<w:spacing w:line="360" w:lineRule="auto"/>' + `
    '<w:jc w:val="center"/>' + `
    '<w:rPr>' + `
      '<w:rFonts w:hint="cs"/>' + `
      '<w:sz w:val="21"/>' + `
      '<w:szCs w:val="21"/>' + `
      '<w:rtl/>' + `
      '<w:lang w:val="en-US"/>' + `
    '</w:rPr>' + `
  '</w:pPr>' + `
  '<w:r w:rsidRPr="00C87AA4">' + `
    '<w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US"/></w:rPr>' + `
    '<w:t>Jonathan Matetzky</w:t>' + `
  '</w:r>' + `
  '<w:r w:rsidRPr="00C87AA4">' + `
    '<w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US"/></w:rPr>' + `
    '<w:tab/><w:t>205949753</w:t>' + `
  '</w:r>' + `
'</w:p>'

$found1 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Jonathan Matetzky*") {
        $p.Range.InsertXML($target1)
        $found1 = $true
        break
    }
}

# ---------------------------------------------------------------------------
# Hunk 2: paragraph containing "26-bit chunks" -- the run is split so that
# "bit" becomes "bytes " (three runs instead of one, same rPr throughout).
# ---------------------------------------------------------------------------
$target2 = '<w:p xmlns:w="' + $wNs + '" xmlns:w14="' + $w14Ns + '" w14:paraId="70421289" w14:textId="3AC92FEB" w:rsidR="00D51C15" w:rsidRDefault="00D51C15" w:rsidP="00D51C15">' + `
  '<w:pPr>' + `
    '<w:pStyle w:val="ListParagraph"/>' + `
    '<w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' + `
    '<w:spacing w:line="360" w:lineRule="auto"/>' + `
    '<w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US"/></w:rPr>' + `
  '</w:pPr>' + `
  '<w:r>' + `
    '<w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US"/></w:rPr>' + `
    '<w:t xml:space="preserve">In case the message file is not in a </w:t>' + `
  '</w:r>' + `
  '<w:r>' + `
    '<w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US"/></w:rPr>' + `
    '<w:t>26-</w:t>' + `
  '</w:r>' + `
  '<w:r>' + `
    '<w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US"/></w:rPr>' + `
    '<w:t xml:space="preserve">bytes </w:t>' + `
  '</w:r>' + `
  '<w:r>' + `
    '<w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US"/></w:rPr>' + `
    '<w:t>chunks</w:t>' + `
  '</w:r>' + `
  '<w:r>' + `
    '<w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US"/></w:rPr>' + `
    '<w:t xml:space="preserve"> format. It won' + [char]0x2019 + 't be able to be read. Same thing wo</w:t>' + `
  '</w:r>' + `
  '<w:r w:rsidR="00E16E8C">' + `
    '<w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US"/></w:rPr>' + `
    '<w:t>u</w:t>' + `
  '</w:r>' + `
  '<w:r>' + `
    '<w:rPr><w:sz w:val="21"/><w:szCs w:val="21"/><w:lang w:val="en-US"/></w:rPr>' + `
    '<w:t>ld happen in case the message would include chars different from ' + [char]0x2018 + '1' + [char]0x2019 + ' or ' + [char]0x2018 + '0' + [char]0x2019 + '.</w:t>' + `
  '</w:r>' + `
'</w:p>'

$found2 = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*26-bit chunks*") {
        $p.Range.InsertXML($target2)
        $found2 = $true
        break
    }
}

Write-Host "Hunk1 applied: $found1; Hunk2 applied: $found2"
